$d = $word.ActiveDocument

$pairs = @(
    @("74×67=", "84×75="),
    @("21×26=", "48×24="),
    @("93×33=", "45×33="),
    @("66×88=", "88×17="),
    @("18×17=", "89×98="),
    @("80×99=", "89×64="),
    @("62×85=", "36×30="),
    @("85×81=", "94×22="),
    @("23×15=", "77×60="),
    @("49×11=", "89×81="),
    @("96×79=", "25×74="),
    @("52×85=", "21×88="),
    @("95×67=", "48×82="),
    @("16×25=", "88×57="),
    @("43×56=", "46×14="),
    @("88×31=", "62×73="),
    @("42×97=", "92×67="),
    @("94×49=", "19×43="),
    @("38×21=", "90×53="),
    @("13×80=", "64×68="),
    @("80×51=", "29×93="),
    @("40×16=", "74×45="),
    @("45×28=", "92×45="),
    @("58×96=", "16×12="),
    @("27×44=", "56×96=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
